# Rename the four duplicated "Лист2 (n)" sheets to their descriptive
# English names (matches the uploaded workbook's new sheet tab labels).
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(3).Name = "normal"
$wb.Worksheets.Item(4).Name = "fnlock"
$wb.Worksheets.Item(5).Name = "modifier"
$wb.Worksheets.Item(6).Name = "media"

# Widen a few columns on the "modifier" / "media" sheets that were resized
# (closest width the host's character-width model can reproduce).
$wsNormal = $wb.Worksheets.Item(3)
$wsFnlock = $wb.Worksheets.Item(4)
$wsModifier = $wb.Worksheets.Item(5)
$wsMedia = $wb.Worksheets.Item(6)

$wsNormal.Columns.Item(6).ColumnWidth = 10.75
$wsFnlock.Columns.Item(6).ColumnWidth = 10.75

$wsModifier.Columns.Item(4).ColumnWidth = 11.25
$wsModifier.Columns.Item(5).ColumnWidth = 11.6
$wsModifier.Columns.Item(10).ColumnWidth = 8.75
$wsModifier.Columns.Item(12).ColumnWidth = 10.25

$wsMedia.Columns.Item(4).ColumnWidth = 11.75
$wsMedia.Columns.Item(10).ColumnWidth = 12.25

# Move each sheet's selection to where the author left it, and leave the
# workbook on "Лист2" (tab 2) as the active/visible sheet.
$wsNormal.Activate()
$null = $wsNormal.Range("E19").Select()

$wsFnlock.Activate()
$null = $wsFnlock.Range("K5").Select()

$wsModifier.Activate()
$null = $wsModifier.Range("F18").Select()

$wsMedia.Activate()
$null = $wsMedia.Range("H18").Select()

$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$null = $ws1.Range("F13").Select()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$null = $ws2.Range("F13").Select()
